$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BPTBfRN")

$newItems = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor",
    "hydrogen"
)

$startRow = 19
for ($i = 0; $i -lt $newItems.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newItems[$i]
    $ws.Cells.Item($row, 2).Value = 1
}

$ws.Activate()
$ws.Range("A25").Select()
